$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

# Update column G ("Recorded By") inline-string values for the affected rows.
# The set of emails per row is unchanged; only the ordering within the
# comma-separated list has been updated to match the upstream source.
$ws.Cells.Item(2, 7).Value = "majorelle.magdy@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg, servinaz@med.asu.edu.eg"
$ws.Cells.Item(10, 7).Value = "esraa.mostafa@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg"
$ws.Cells.Item(15, 7).Value = "hana.amr@med.asu.edu.eg, nancy.abdelshafy@med.asu.edu.eg"
$ws.Cells.Item(18, 7).Value = "Remon.Matta@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg"
$ws.Cells.Item(19, 7).Value = "naema.gomaa@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg"
$ws.Cells.Item(20, 7).Value = "aya.hanafy@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, Wafaa.ebida@med.asu.edu.eg, marinasorial@med.asu.edu.eg"
$ws.Cells.Item(24, 7).Value = "majorelle.magdy@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg, servinaz@med.asu.edu.eg"
$ws.Cells.Item(25, 7).Value = "alshimaa.atef@med.asu.edu.egm, backup@backdoor.com, manar.montaser@med.asu.edu.eg"
$ws.Cells.Item(32, 7).Value = "esraa.mostafa@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg"
$ws.Cells.Item(37, 7).Value = "hana.amr@med.asu.edu.eg, nancy.abdelshafy@med.asu.edu.eg"
$ws.Cells.Item(40, 7).Value = "Remon.Matta@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg"
$ws.Cells.Item(41, 7).Value = "naema.gomaa@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg"
$ws.Cells.Item(42, 7).Value = "aya.hanafy@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, Wafaa.ebida@med.asu.edu.eg, marinasorial@med.asu.edu.eg"
$ws.Cells.Item(46, 7).Value = "nourhan.mahmoud@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg"
$ws.Cells.Item(54, 7).Value = "arwaelsayed03@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, maimustafa@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, merna.said@med.asu.edu.eg"
$ws.Cells.Item(62, 7).Value = "yassmen.ahmed@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg"
$ws.Cells.Item(63, 7).Value = "ola.m.abdelfattah@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg"
$ws.Cells.Item(64, 7).Value = "youstina.magdy@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg"
$ws.Cells.Item(68, 7).Value = "nourhan.mahmoud@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg"
$ws.Cells.Item(76, 7).Value = "arwaelsayed03@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, maimustafa@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, merna.said@med.asu.edu.eg"
$ws.Cells.Item(84, 7).Value = "yassmen.ahmed@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg"
$ws.Cells.Item(85, 7).Value = "ola.m.abdelfattah@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg"
$ws.Cells.Item(86, 7).Value = "youstina.magdy@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg"
$ws.Cells.Item(90, 7).Value = "shaimaa.ahmed@med.asu.edu.eg, manar.montaser@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg"
$ws.Cells.Item(96, 7).Value = "aml.awwad@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg, norhan.mohamed@med.asu.edu.eg"
$ws.Cells.Item(98, 7).Value = "arwaelsayed03@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, maimustafa@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, merna.said@med.asu.edu.eg"
$ws.Cells.Item(106, 7).Value = "wafaa.ebida@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg"
$ws.Cells.Item(107, 7).Value = "aya.hanafy@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg"
$ws.Cells.Item(108, 7).Value = "aya.hanafy@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, Wafaa.ebida@med.asu.edu.eg, marinasorial@med.asu.edu.eg"
$ws.Cells.Item(112, 7).Value = "shaimaa.ahmed@med.asu.edu.eg, manar.montaser@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg"
$ws.Cells.Item(118, 7).Value = "aml.awwad@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg, norhan.mohamed@med.asu.edu.eg"
$ws.Cells.Item(120, 7).Value = "arwaelsayed03@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, maimustafa@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, merna.said@med.asu.edu.eg"
$ws.Cells.Item(128, 7).Value = "wafaa.ebida@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg"
$ws.Cells.Item(129, 7).Value = "aya.hanafy@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg"
$ws.Cells.Item(130, 7).Value = "aya.hanafy@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, Wafaa.ebida@med.asu.edu.eg, marinasorial@med.asu.edu.eg"
$ws.Cells.Item(134, 7).Value = "majorelle.magdy@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg"
$ws.Cells.Item(142, 7).Value = "esraa.mostafa@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, merna.said@med.asu.edu.eg"
$ws.Cells.Item(150, 7).Value = "naema.gomaa@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg"
$ws.Cells.Item(151, 7).Value = "yassmen.ahmed@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg"
$ws.Cells.Item(156, 7).Value = "majorelle.magdy@med.asu.edu.eg, alshimaa.atef@med.asu.edu.egm, manar.montaser@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg"
$ws.Cells.Item(164, 7).Value = "esraa.mostafa@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, merna.said@med.asu.edu.eg"
$ws.Cells.Item(172, 7).Value = "naema.gomaa@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg"
$ws.Cells.Item(173, 7).Value = "yassmen.ahmed@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg"
